# ex9.1.2(Linear)__M_Stationarygenerator_alpha_non_zero.xlsx
# "expermits todos no convexos menos el 5to"
# Rewrite the generated experiment values (restrictions, modified point,
# bf/BF/alpha vectors) with a new non-convex generator run.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Restricciones_del_follower (sheet3): rows 2-5, columns A-F
# Force text storage for the numeric-looking values (they are stored as
# plain text in the workbook, not as numbers) by pre-setting the
# NumberFormat to "@" on the whole block before writing.
# ---------------------------------------------------------------------
$wsFollower = $wb.Worksheets.Item("Restricciones_del_follower")

$wsFollower.Range("A2:B2").NumberFormat = "@"
$wsFollower.Range("D2:F2").NumberFormat = "@"
$wsFollower.Range("B3").NumberFormat = "@"
$wsFollower.Range("D3:F3").NumberFormat = "@"
$wsFollower.Range("B4").NumberFormat = "@"
$wsFollower.Range("D4:E4").NumberFormat = "@"
$wsFollower.Range("B5").NumberFormat = "@"
$wsFollower.Range("D5:F5").NumberFormat = "@"

$wsFollower.Range("A2").Value = "0"
$wsFollower.Range("B2").Value = "0"
$wsFollower.Range("D2").Value = "0.13"
$wsFollower.Range("E2").Value = "2.9"
$wsFollower.Range("F2").Value = "0"

$wsFollower.Range("A3").Value = "-24.29684212488997 - x + 3.776005232905586y"
$wsFollower.Range("B3").Value = "21.29684212488997"
$wsFollower.Range("D3").Value = "0.6"
$wsFollower.Range("E3").Value = "1.7000000000000002"
$wsFollower.Range("F3").Value = "5.1"

$wsFollower.Range("A4").Value = "-0.4547970843733591 + x - 0.7012674429349861y"
$wsFollower.Range("B4").Value = "-11.545202915626641"
$wsFollower.Range("D4").Value = "0.62"
$wsFollower.Range("E4").Value = "9.0"

$wsFollower.Range("A5").Value = "-50.53871002686802 + 4x + 3.1836906865674557y"
$wsFollower.Range("B5").Value = "38.02871002686802"
$wsFollower.Range("D5").Value = "0.32"
$wsFollower.Range("E5").Value = "8.100000000000001"
$wsFollower.Range("F5").Value = "4.3"

# ---------------------------------------------------------------------
# Punto_modificado (sheet4): new modified point (x, y)
# ---------------------------------------------------------------------
$wsPunto = $wb.Worksheets.Item("Punto_modificado")
$wsPunto.Range("A2:B2").NumberFormat = "@"
$wsPunto.Range("A2").Value = "6.1"
$wsPunto.Range("B2").Value = "8.05"

# ---------------------------------------------------------------------
# Vector_bf (sheet5) / Vector_BF (sheet6)
# NOTE: the two sheet names only differ by case ("Vector_bf" vs
# "Vector_BF"); worksheet name lookups are case-insensitive, so both
# names would resolve to the same sheet object. Use the stable 1-based
# tab index instead to address each sheet unambiguously.
# ---------------------------------------------------------------------
$wsBf = $wb.Worksheets.Item(5)
$wsBf.Range("A2").NumberFormat = "@"
$wsBf.Range("A2").Value = "-3.849598344825246"

$wsBF = $wb.Worksheets.Item(6)
$wsBF.Range("A2:A3").NumberFormat = "@"
$wsBF.Range("A2").Value = "-38.7"
$wsBF.Range("A3").Value = "-22.895696470721017"

# ---------------------------------------------------------------------
# Vector_Alpha (sheet7): A2 stays a genuine number
# ---------------------------------------------------------------------
$wsAlpha = $wb.Worksheets.Item("Vector_Alpha")
$wsAlpha.Range("A2").Value = 1.350633721467493
